$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update N (latitude) and O (longitude) columns for rows whose M column
# holds one of the recognized place names (Trai / Cau Cay Bang / Vung / Bo Coc).

$ws.Range("N2").Style = "Normal"
$ws.Range("N2").HorizontalAlignment = -4131
$ws.Range("N2").Value2 = 21.021014219286901
$ws.Range("O2").Style = "Normal"
$ws.Range("O2").Value2 = 105.666575814257

$ws.Range("N4").Style = "Normal"
$ws.Range("N4").HorizontalAlignment = -4131
$ws.Range("N4").Value2 = 21.021014219286901
$ws.Range("O4").Style = "Normal"
$ws.Range("O4").Value2 = 105.666575814257

$ws.Range("N6").Style = "Normal"
$ws.Range("N6").HorizontalAlignment = -4131
$ws.Range("N6").Value2 = 21.021014219286901
$ws.Range("O6").Style = "Normal"
$ws.Range("O6").Value2 = 105.666575814257

$ws.Range("N8").Style = "Normal"
$ws.Range("N8").HorizontalAlignment = -4131
$ws.Range("N8").Value2 = 21.021014219286901
$ws.Range("O8").Style = "Normal"
$ws.Range("O8").Value2 = 105.666575814257

$ws.Range("N10").Style = "Normal"
$ws.Range("N10").HorizontalAlignment = -4131
$ws.Range("N10").Value2 = 21.0129386402952
$ws.Range("O10").Style = "Normal"
$ws.Range("O10").Value2 = 105.65060486244001

$ws.Range("N13").Style = "Normal"
$ws.Range("N13").HorizontalAlignment = -4131
$ws.Range("N13").Value2 = 21.016309601299501
$ws.Range("O13").Style = "Normal"
$ws.Range("O13").Value2 = 105.65990582900299

$ws.Range("N14").Style = "Normal"
$ws.Range("N14").HorizontalAlignment = -4131
$ws.Range("N14").Value2 = 21.016309601299501
$ws.Range("O14").Style = "Normal"
$ws.Range("O14").Value2 = 105.65990582900299

$ws.Range("N30").Style = "Normal"
$ws.Range("N30").HorizontalAlignment = -4131
$ws.Range("N30").Value2 = 21.0129386402952
$ws.Range("O30").Style = "Normal"
$ws.Range("O30").Value2 = 105.65060486244001

$ws.Range("N37").Style = "Normal"
$ws.Range("N37").HorizontalAlignment = -4131
$ws.Range("N37").Value2 = 21.017046103716801
$ws.Range("O37").Style = "Normal"
$ws.Range("O37").Value2 = 105.655111203495

$ws.Range("N40").Style = "Normal"
$ws.Range("N40").HorizontalAlignment = -4131
$ws.Range("N40").Value2 = 21.021014219286901
$ws.Range("O40").Style = "Normal"
$ws.Range("O40").Value2 = 105.666575814257

$ws.Range("N43").Style = "Normal"
$ws.Range("N43").HorizontalAlignment = -4131
$ws.Range("N43").Value2 = 21.016309601299501
$ws.Range("O43").Style = "Normal"
$ws.Range("O43").Value2 = 105.65990582900299

$ws.Range("N46").Style = "Normal"
$ws.Range("N46").HorizontalAlignment = -4131
$ws.Range("N46").Value2 = 21.0129386402952
$ws.Range("O46").Style = "Normal"
$ws.Range("O46").Value2 = 105.65060486244001

$ws.Range("N47").Style = "Normal"
$ws.Range("N47").HorizontalAlignment = -4131
$ws.Range("N47").Value2 = 21.0129386402952
$ws.Range("O47").Style = "Normal"
$ws.Range("O47").Value2 = 105.65060486244001

$ws.Range("N50").Style = "Normal"
$ws.Range("N50").HorizontalAlignment = -4131
$ws.Range("N50").Value2 = 21.0129386402952
$ws.Range("O50").Style = "Normal"
$ws.Range("O50").Value2 = 105.65060486244001

$ws.Range("N54").Style = "Normal"
$ws.Range("N54").HorizontalAlignment = -4131
$ws.Range("N54").Value2 = 21.0129386402952
$ws.Range("O54").Style = "Normal"
$ws.Range("O54").Value2 = 105.65060486244001

$ws.Range("N63").Style = "Normal"
$ws.Range("N63").HorizontalAlignment = -4131
$ws.Range("N63").Value2 = 21.021014219286901
$ws.Range("O63").Style = "Normal"
$ws.Range("O63").Value2 = 105.666575814257

$ws.Range("N64").Style = "Normal"
$ws.Range("N64").HorizontalAlignment = -4131
$ws.Range("N64").Value2 = 21.0129386402952
$ws.Range("O64").Style = "Normal"
$ws.Range("O64").Value2 = 105.65060486244001

$ws.Range("N74").Style = "Normal"
$ws.Range("N74").HorizontalAlignment = -4131
$ws.Range("N74").Value2 = 21.021014219286901
$ws.Range("O74").Style = "Normal"
$ws.Range("O74").Value2 = 105.666575814257

$ws.Range("N75").Style = "Normal"
$ws.Range("N75").HorizontalAlignment = -4131
$ws.Range("N75").Value2 = 21.017046103716801
$ws.Range("O75").Style = "Normal"
$ws.Range("O75").Value2 = 105.655111203495

$ws.Range("N84").Style = "Normal"
$ws.Range("N84").HorizontalAlignment = -4131
$ws.Range("N84").Value2 = 21.021014219286901
$ws.Range("O84").Style = "Normal"
$ws.Range("O84").Value2 = 105.666575814257

$ws.Range("N85").Style = "Normal"
$ws.Range("N85").HorizontalAlignment = -4131
$ws.Range("N85").Value2 = 21.0129386402952
$ws.Range("O85").Style = "Normal"
$ws.Range("O85").Value2 = 105.65060486244001

$ws.Range("N92").Style = "Normal"
$ws.Range("N92").HorizontalAlignment = -4131
$ws.Range("N92").Value2 = 21.017046103716801
$ws.Range("O92").Style = "Normal"
$ws.Range("O92").Value2 = 105.655111203495

$ws.Range("N93").Style = "Normal"
$ws.Range("N93").HorizontalAlignment = -4131
$ws.Range("N93").Value2 = 21.021014219286901
$ws.Range("O93").Style = "Normal"
$ws.Range("O93").Value2 = 105.666575814257

$ws.Range("N95").Style = "Normal"
$ws.Range("N95").HorizontalAlignment = -4131
$ws.Range("N95").Value2 = 21.0129386402952
$ws.Range("O95").Style = "Normal"
$ws.Range("O95").Value2 = 105.65060486244001

$ws.Range("N96").Style = "Normal"
$ws.Range("N96").HorizontalAlignment = -4131
$ws.Range("N96").Value2 = 21.0129386402952
$ws.Range("O96").Style = "Normal"
$ws.Range("O96").Value2 = 105.65060486244001

$ws.Range("N98").Style = "Normal"
$ws.Range("N98").HorizontalAlignment = -4131
$ws.Range("N98").Value2 = 21.0129386402952
$ws.Range("O98").Style = "Normal"
$ws.Range("O98").Value2 = 105.65060486244001

$ws.Range("N101").Style = "Normal"
$ws.Range("N101").HorizontalAlignment = -4131
$ws.Range("N101").Value2 = 21.0129386402952
$ws.Range("O101").Style = "Normal"
$ws.Range("O101").Value2 = 105.65060486244001

$ws.Range("N107").Style = "Normal"
$ws.Range("N107").HorizontalAlignment = -4131
$ws.Range("N107").Value2 = 21.0129386402952
$ws.Range("O107").Style = "Normal"
$ws.Range("O107").Value2 = 105.65060486244001

$ws.Range("N113").Style = "Normal"
$ws.Range("N113").HorizontalAlignment = -4131
$ws.Range("N113").Value2 = 21.0129386402952
$ws.Range("O113").Style = "Normal"
$ws.Range("O113").Value2 = 105.65060486244001

$ws.Range("N115").Style = "Normal"
$ws.Range("N115").HorizontalAlignment = -4131
$ws.Range("N115").Value2 = 21.0129386402952
$ws.Range("O115").Style = "Normal"
$ws.Range("O115").Value2 = 105.65060486244001

$ws.Range("N116").Style = "Normal"
$ws.Range("N116").HorizontalAlignment = -4131
$ws.Range("N116").Value2 = 21.017046103716801
$ws.Range("O116").Style = "Normal"
$ws.Range("O116").Value2 = 105.655111203495

$ws.Range("N117").Style = "Normal"
$ws.Range("N117").HorizontalAlignment = -4131
$ws.Range("N117").Value2 = 21.017046103716801
$ws.Range("O117").Style = "Normal"
$ws.Range("O117").Value2 = 105.655111203495

$ws.Range("N120").Style = "Normal"
$ws.Range("N120").HorizontalAlignment = -4131
$ws.Range("N120").Value2 = 21.0129386402952
$ws.Range("O120").Style = "Normal"
$ws.Range("O120").Value2 = 105.65060486244001

$ws.Range("N121").Style = "Normal"
$ws.Range("N121").HorizontalAlignment = -4131
$ws.Range("N121").Value2 = 21.017046103716801
$ws.Range("O121").Style = "Normal"
$ws.Range("O121").Value2 = 105.655111203495

$ws.Range("N126").Style = "Normal"
$ws.Range("N126").HorizontalAlignment = -4131
$ws.Range("N126").Value2 = 21.017046103716801
$ws.Range("O126").Style = "Normal"
$ws.Range("O126").Value2 = 105.655111203495

$ws.Range("N129").Style = "Normal"
$ws.Range("N129").HorizontalAlignment = -4131
$ws.Range("N129").Value2 = 21.0129386402952
$ws.Range("O129").Style = "Normal"
$ws.Range("O129").Value2 = 105.65060486244001

$ws.Range("N130").Style = "Normal"
$ws.Range("N130").HorizontalAlignment = -4131
$ws.Range("N130").Value2 = 21.0129386402952
$ws.Range("O130").Style = "Normal"
$ws.Range("O130").Value2 = 105.65060486244001

$ws.Range("N132").Style = "Normal"
$ws.Range("N132").HorizontalAlignment = -4131
$ws.Range("N132").Value2 = 21.0129386402952
$ws.Range("O132").Style = "Normal"
$ws.Range("O132").Value2 = 105.65060486244001

$ws.Range("N133").Style = "Normal"
$ws.Range("N133").HorizontalAlignment = -4131
$ws.Range("N133").Value2 = 21.017046103716801
$ws.Range("O133").Style = "Normal"
$ws.Range("O133").Value2 = 105.655111203495

$ws.Range("N134").Style = "Normal"
$ws.Range("N134").HorizontalAlignment = -4131
$ws.Range("N134").Value2 = 21.017046103716801
$ws.Range("O134").Style = "Normal"
$ws.Range("O134").Value2 = 105.655111203495

$ws.Range("N137").Style = "Normal"
$ws.Range("N137").HorizontalAlignment = -4131
$ws.Range("N137").Value2 = 21.0129386402952
$ws.Range("O137").Style = "Normal"
$ws.Range("O137").Value2 = 105.65060486244001

$ws.Range("N140").Style = "Normal"
$ws.Range("N140").HorizontalAlignment = -4131
$ws.Range("N140").Value2 = 21.017046103716801
$ws.Range("O140").Style = "Normal"
$ws.Range("O140").Value2 = 105.655111203495

$ws.Range("N142").Style = "Normal"
$ws.Range("N142").HorizontalAlignment = -4131
$ws.Range("N142").Value2 = 21.017046103716801
$ws.Range("O142").Style = "Normal"
$ws.Range("O142").Value2 = 105.655111203495

$ws.Range("N143").Style = "Normal"
$ws.Range("N143").HorizontalAlignment = -4131
$ws.Range("N143").Value2 = 21.017046103716801
$ws.Range("O143").Style = "Normal"
$ws.Range("O143").Value2 = 105.655111203495

$ws.Range("N160").Style = "Normal"
$ws.Range("N160").HorizontalAlignment = -4131
$ws.Range("N160").Value2 = 21.0129386402952
$ws.Range("O160").Style = "Normal"
$ws.Range("O160").Value2 = 105.65060486244001

$ws.Range("N161").Style = "Normal"
$ws.Range("N161").HorizontalAlignment = -4131
$ws.Range("N161").Value2 = 21.017046103716801
$ws.Range("O161").Style = "Normal"
$ws.Range("O161").Value2 = 105.655111203495

$ws.Range("N164").Style = "Normal"
$ws.Range("N164").HorizontalAlignment = -4131
$ws.Range("N164").Value2 = 21.017046103716801
$ws.Range("O164").Style = "Normal"
$ws.Range("O164").Value2 = 105.655111203495

$ws.Range("N168").Style = "Normal"
$ws.Range("N168").HorizontalAlignment = -4131
$ws.Range("N168").Value2 = 21.021014219286901
$ws.Range("O168").Style = "Normal"
$ws.Range("O168").Value2 = 105.666575814257

$ws.Range("N169").Style = "Normal"
$ws.Range("N169").HorizontalAlignment = -4131
$ws.Range("N169").Value2 = 21.021014219286901
$ws.Range("O169").Style = "Normal"
$ws.Range("O169").Value2 = 105.666575814257

$ws.Range("N171").Style = "Normal"
$ws.Range("N171").HorizontalAlignment = -4131
$ws.Range("N171").Value2 = 21.0129386402952
$ws.Range("O171").Style = "Normal"
$ws.Range("O171").Value2 = 105.65060486244001

$ws.Range("N172").Style = "Normal"
$ws.Range("N172").HorizontalAlignment = -4131
$ws.Range("N172").Value2 = 21.0129386402952
$ws.Range("O172").Style = "Normal"
$ws.Range("O172").Value2 = 105.65060486244001

$ws.Range("N179").Style = "Normal"
$ws.Range("N179").HorizontalAlignment = -4131
$ws.Range("N179").Value2 = 21.0129386402952
$ws.Range("O179").Style = "Normal"
$ws.Range("O179").Value2 = 105.65060486244001

$ws.Range("N188").Style = "Normal"
$ws.Range("N188").HorizontalAlignment = -4131
$ws.Range("N188").Value2 = 21.017046103716801
$ws.Range("O188").Style = "Normal"
$ws.Range("O188").Value2 = 105.655111203495

$ws.Range("N189").Style = "Normal"
$ws.Range("N189").HorizontalAlignment = -4131
$ws.Range("N189").Value2 = 21.0129386402952
$ws.Range("O189").Style = "Normal"
$ws.Range("O189").Value2 = 105.65060486244001

$ws.Range("N203").Style = "Normal"
$ws.Range("N203").HorizontalAlignment = -4131
$ws.Range("N203").Value2 = 21.017046103716801
$ws.Range("O203").Style = "Normal"
$ws.Range("O203").Value2 = 105.655111203495

$ws.Range("N211").Style = "Normal"
$ws.Range("N211").HorizontalAlignment = -4131
$ws.Range("N211").Value2 = 21.017046103716801
$ws.Range("O211").Style = "Normal"
$ws.Range("O211").Value2 = 105.655111203495

$ws.Range("N215").Style = "Normal"
$ws.Range("N215").HorizontalAlignment = -4131
$ws.Range("N215").Value2 = 21.0129386402952
$ws.Range("O215").Style = "Normal"
$ws.Range("O215").Value2 = 105.65060486244001

$ws.Range("N233").Style = "Normal"
$ws.Range("N233").HorizontalAlignment = -4131
$ws.Range("N233").Value2 = 21.017046103716801
$ws.Range("O233").Style = "Normal"
$ws.Range("O233").Value2 = 105.655111203495

$ws.Range("N237").Style = "Normal"
$ws.Range("N237").HorizontalAlignment = -4131
$ws.Range("N237").Value2 = 21.021014219286901
$ws.Range("O237").Style = "Normal"
$ws.Range("O237").Value2 = 105.666575814257

$ws.Range("N370").Style = "Normal"
$ws.Range("N370").HorizontalAlignment = -4131
$ws.Range("N370").Value2 = 21.0129386402952
$ws.Range("O370").Style = "Normal"
$ws.Range("O370").Value2 = 105.65060486244001

# Update the sheet view: selection moves to N462:O462 (pane top-left row follows automatically).
$ws.Range("N462:O462").Select()
